$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove one row of data (old row 4 held the password that no longer
#    appears in the refreshed list) and shrink the used range from 6 to 5
#    rows before rewriting the sheet with the new data set.
# ---------------------------------------------------------------------------
$ws.Rows("4").Delete()

# ---------------------------------------------------------------------------
# 2. Headers: A keeps the file-path column (renamed), B keeps the password
#    column (same caption), C is a brand-new "file name" column.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "CAMINHO ARQUIVO"
$ws.Range("B1").Value = "SENHA ARQUIVO"
$ws.Range("C1").Value = "NOME ARQUIVO"

# ---------------------------------------------------------------------------
# 3. Data rows: new file paths in column A, same passwords in column B.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\C1852538_SPI.pdf"
$ws.Range("B2").Value = 12345

$ws.Range("A3").Value = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\C1852539_SPI.pdf"
$ws.Range("B3").Value = 54321

$ws.Range("A4").Value = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\C1852554_SPI.pdf"
$ws.Range("B4").Value = 32145

$ws.Range("A5").Value = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\C1852789_SPI.pdf"
$ws.Range("B5").Value = 32154

# ---------------------------------------------------------------------------
# 4. Column C: formula extracting the last 16 characters (the file name)
#    from the full path in column A.
# ---------------------------------------------------------------------------
$ws.Range("C2").Formula = "=RIGHT(A2,16)"
$ws.Range("C3:C5").Formula = "=RIGHT(A3,16)"

# ---------------------------------------------------------------------------
# 5. Column widths - new column C mirrors column B's width.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 18.28515625

# ---------------------------------------------------------------------------
# 6. Formatting: bold header row, thick box border around header cells,
#    thin box border around data cells, and alignment.
# ---------------------------------------------------------------------------
$ws.Rows("1").RowHeight = 31.5

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 4

$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").VerticalAlignment = -4108

$dataRange = $ws.Range("A2:C5")
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# ---------------------------------------------------------------------------
# 7. Misc workbook/view bookkeeping to mirror the authored change.
# ---------------------------------------------------------------------------
$ws.Range("F13").Select()
